# Update the "取得日時" (retrieved datetime) timestamps in column A of the
# "ランサーズ" sheet for rows 2 through 14 from 2025-09-22 06:27:57 to
# 2025-09-22 06:37:29, matching the commit "Append: 2025-09-22 06:37 JST".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-09-22 06:27:57"
$newValue = "2025-09-22 06:37:29"

for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
